# ---------------------------------------------------------------------------
# "The version has a working Add method in the CustomList class"
#
# Applies two edits:
#  1. Splits the "Capacity property ... size of my " run so the (hidden)
#     _GoBack bookmark now sits right after "...size of" and before " my ".
#  2. Rewrites the "Add Method" scratch notes at the end of the document:
#       - "1)Create a list" -> full sentence about declaring/instantiating.
#       - adds "But then we will need to:" / "1) Get list count" /
#         "2) Get maximum values for Array" / the long "3) If the Count..."
#         paragraph, moving the old trailing _GoBack bookmark away from the
#         final "2)" paragraph (it now lives where step 1 put it).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Split-RunAt($pos) {
    # Inserting (and immediately removing) a bookmark at a zero-length
    # range forces Word to split the run at that character offset without
    # altering any text, so adjacent runs keep identical formatting but
    # remain distinct <w:r> elements (mirrors how the diff shows two
    # same-rPr runs sitting side by side).
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("__SplitHelper__", $bmRange)
    $d.Bookmarks("__SplitHelper__").Delete()
}

# ---------------------------------------------------------------------------
# Step 2 (done first on purpose): rewrite "1)Create a list" and build the
# four new paragraphs that follow it. All the paragraph breaks are created
# up front -- while still in the original (inherited, non-explicit) bold
# run -- so that only the paragraphs that really need an explicit bold
# paragraph-mark (<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>) get one.
# ---------------------------------------------------------------------------

$p1 = $d.Content
$p1.Find.Execute("1)Create a list", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$para1 = $p1.Paragraphs(1)
$p1.Text = "First need to make sure a list is created by declaring the list and instantiating the list."

$para1.Range.InsertParagraphAfter()
$para2 = $para1.Next()
$para2.Range.InsertParagraphAfter()
$para3 = $para2.Next()
$para4 = $para3.Next()          # the original "2)" paragraph
$para4.Range.InsertParagraphAfter()
$para5 = $d.Paragraphs.Last

# Paragraph 2: "But then we will need to:"
$para2.Range.Font.Bold = 1
$para2.Range.InsertAfter("But then we will need to:")

# Paragraph 3: "1)" followed by " Get list count" as two separate runs.
$para3.Range.Font.Bold = 1
$p3Start = $para3.Range.Start
$para3.Range.InsertAfter("1) Get list count")
Split-RunAt ($p3Start + 2)

# Paragraph 4: "2)" (already present) followed by " Get maximum values for
# Array" as a second run.
$para4.Range.Font.Bold = 1
$p4Start = $para4.Range.Start
$para4.Range.InsertAfter(" Get maximum values for Array")
Split-RunAt ($p4Start + 2)

# Paragraph 5: the long "3) If the Count..." sentence -- bold is inherited
# from paragraph 4 without an explicit Font.Bold write, so no paragraph
# mark formatting gets recorded for this paragraph.
$fullText = "3) If the Count of the Array would make it equal to the Max Capacity of the array then we need to create a new array which will double the array size.  It will copy all of the values for the old array, and add them in the appropriate indexes of the new array.  "
$para5.Range.InsertAfter($fullText)

# ---------------------------------------------------------------------------
# Step 1 (done last so nothing else re-homes the hidden bookmark afterwards):
# move the _GoBack bookmark into the Capacity-property paragraph, right
# between "...size of" and " my ".
# ---------------------------------------------------------------------------

$findRange = $d.Content
$findRange.Find.Execute("so that I can publicly see the size of", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$splitPos = $findRange.End

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

Write-Host "done"
